$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting (style index) of row 2 data cells (B2:C2) onto row 3 (B3:C3) ---
# This reuses the existing cellXfs entry (fontId=2, wrapText) instead of minting a new one.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B3:C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 3 values ---
$ws.Range("A3").Value = "FTP"
$ws.Range("B3").Value = "example"

$exampleText = @'
>ftp 1?2.18.255.20
Connected to 172.18.255.20.
200-QTCP at S657274B.
200 Connection will close if idle more than 5 minutes. 
User (172.18.255.20:(none)): ZCSERUICE
331 Enter password. 
Password: xxxx
230 ZCSERUICE logged on. 
ftp> CD YMYLES
250 ‘YMYLES’ is current library. 
ftp> PWD
257 ‘YMYLES’ is current library..
ftp> LCD D:\
Local directory now D:\. 
ftp> PUT D:\20161205.TXT
200 PORT subcommand request successful.
150-Characters that are not recognized in file specification are changed to #.
150 Sending file to member TXT in file 110161205 in library YMYLES.
226 File transfer completed successfully.
ftp: 1237 bytes sent in 0.27Seconds 4.67Mbytes/sec.
ftp> quit
'@
$ws.Range("C3").Value = $exampleText

# --- Rich-text run formatting inside C3, matching the captured FTP session transcript ---
# (red = user-typed ftp commands, black = server/system output, "23" kept colour-less like the source)
$c = $ws.Range("C3").Characters(2, 17)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(19, 136)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(155, 9)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(164, 20)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(184, 1)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(185, 11)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(196, 5)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(201, 2)
$c.Font.Size = 8
$c.Font.ColorIndex = -4105

$c = $ws.Range("C3").Characters(203, 29)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(232, 9)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(241, 40)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(281, 3)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(284, 39)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(323, 8)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(331, 32)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(363, 19)
$c.Font.Size = 8
$c.Font.Color = 255

$c = $ws.Range("C3").Characters(382, 287)
$c.Font.Size = 8
$c.Font.Color = 0

$c = $ws.Range("C3").Characters(669, 4)
$c.Font.Size = 8
$c.Font.Color = 255

# --- Row height (the cell wraps a long multi-line transcript) ---
$ws.Rows.Item(3).RowHeight = 237

# --- Selection / view state ---
$ws.Range("C3").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2

